# Generate Report for Handoff
#
# Mark the "3db0f42a-..." file (row 3 on each sheet) as "Ready for handoff"
# instead of "Handed back: in sync with en-US", and refresh its handoff
# timestamps on the Overview sheet as well as the per-language (zh-cn /
# de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the 3db0f42a-...md entry ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-45-13 02:45:14"

# --- zh-cn detail sheet: row 3 is the 3db0f42a-...md entry ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-13 02:45:10"
$wsZhCn.Range("G3").Value = "3db0f42a-77e3-4615-8e75-beae3b4e46b9.09797310f88068644f599538f356718f9bef3e45.zh-cn.xlf"

# --- de-de detail sheet: row 3 is the 3db0f42a-...md entry ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-13 02:45:14"
$wsDeDe.Range("G3").Value = "3db0f42a-77e3-4615-8e75-beae3b4e46b9.09797310f88068644f599538f356718f9bef3e45.de-de.xlf"
